$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeRef, $val) {
    $cell = $ws.Range($rangeRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue "D2" "63.726.22"
$ws.Range("E2").Value = "  -1.72%  "
Set-TextValue "D3" "3.390.59"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "569.78"
$ws.Range("E5").Value = "  -0.78%  "
Set-TextValue "D6" "161.33"
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("E7").Value = "  +0.03%  "
Set-TextValue "D8" "3.389.74"
$ws.Range("E8").Value = "  -1.76%  "
Set-TextValue "D9" "0.544"
$ws.Range("E9").Value = "  -5.36%  "
$ws.Range("E10").Value = "  +1.43%  "
Set-TextValue "D11" "0.118"
$ws.Range("E11").Value = "  -2.54%  "
Set-TextValue "D12" "0.418"
$ws.Range("E12").Value = "  -4.77%  "
Set-TextValue "D13" "3.975.19"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("E14").Value = "  +0.67%  "
Set-TextValue "D15" "26.83"
$ws.Range("E15").Value = "  -2.85%  "
$ws.Range("E16").Value = "  -0.96%  "
Set-TextValue "D17" "63.782.62"
$ws.Range("E17").Value = "  -1.72%  "
Set-TextValue "D18" "3.379.00"
$ws.Range("E18").Value = "  -2.32%  "
Set-TextValue "D19" "6.08"
Set-TextValue "D20" "13.48"
$ws.Range("E20").Value = "  -1.75%  "
Set-TextValue "D21" "375.58"
$ws.Range("E21").Value = "  -0.28%  "
Set-TextValue "D23" "1.00"
$ws.Range("E23").Value = "  +0.06%  "
Set-TextValue "D24" "69.95"
$ws.Range("E24").Value = "  -3.21%  "
Set-TextValue "D25" "0.509"
$ws.Range("E25").Value = "  -4.74%  "
Set-TextValue "D26" "0.0000114"
$ws.Range("E26").Value = "  -5.63%  "
$ws.Range("E27").Value = "  -4.09%  "
Set-TextValue "D28" "0.178"
Set-TextValue "D29" "0.998"
$ws.Range("E29").Value = "  -0.12%  "
Set-TextValue "D30" "6.05"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  -4.02%  "
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("E33").Value = "  +0.04%  "
Set-TextValue "D34" "22.70"
$ws.Range("E34").Value = "  -1.90%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D36" "159.49"
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D37" "1.47"
$ws.Range("E37").Value = "  -6.18%  "
Set-TextValue "D38" "0.857"
$ws.Range("E38").Value = "  +8.83%  "
Set-TextValue "D39" "1.79"
$ws.Range("E39").Value = "  -4.24%  "
Set-TextValue "D40" "0.0719"
$ws.Range("E40").Value = "  -3.97%  "
Set-TextValue "D41" "25.65"
$ws.Range("E41").Value = "  -2.11%  "
Set-TextValue "D42" "42.61"
$ws.Range("E42").Value = "  -0.96%  "
Set-TextValue "D43" "2.732.33"
$ws.Range("E43").Value = "  -5.71%  "
Set-TextValue "D44" "26.00"
$ws.Range("E44").Value = "  -0.78%  "
Set-TextValue "D45" "6.38"
$ws.Range("E45").Value = "  -0.89%  "
Set-TextValue "D46" "4.34"
$ws.Range("E46").Value = "  -3.90%  "
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("E48").Value = "  +1.21%  "
Set-TextValue "D49" "327.33"
$ws.Range("E49").Value = "  +1.84%  "
Set-TextValue "D50" "1.03"
$ws.Range("E50").Value = "  -5.02%  "
$ws.Range("E51").Value = "  -1.97%  "
